$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "CCAA nombre"
$ws.Range("B1").Value = "Siglas"
$ws.Range("C1").Value = "CCAA código"
$ws.Range("D1").Value = "Diputados"
$ws.Range("E1").Value = "Año"
$ws.Range("F1").Value = "Votos"

# Row 2 - concept identifiers
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:siglas"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:diputados"
$ws.Range("E2").Value = "sdmx-dimension:refPeriod"
$ws.Range("F2").Value = "iaest-measure:votos"

# Row 3 - dim / medida
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"

# Row 4 - types
$ws.Range("A4").Value = "URI-Comunidad"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "xsd:date"
$ws.Range("F4").Value = "xsd:int"

# Row 5 - mapping file references
$ws.Range("A5:B5").Delete()
$ws.Range("E5").Value = "mapping-ano.xlsx"

# Match the styling (s="1") used by the other data cells
$ws.Range("A4").Copy()
$ws.Range("E5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
